$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Set the print area to cover both recipe tables (A1:J27) ---
# Use the Names collection so the sheet-qualified reference keeps its
# quoting (matches how Excel itself serializes '500mL'!$A$1:$J$27).
$found = $false
foreach ($n in $wb.Names) {
    if ($n.Name -eq "500mL!Print_Area") {
        $n.RefersTo = "='500mL'!`$A`$1:`$J`$27"
        $found = $true
    }
}
if (-not $found) {
    $ws.PageSetup.PrintArea = '$A$1:$J$27'
}

# --- Scale the printout to 72% (keeps existing landscape orientation) ---
$ws.PageSetup.Zoom = 72

# --- Add the missing bordered (blank) cell above the second table ---
$a11 = $ws.Range("A11")
$a11.Borders.LineStyle = 1
$a11.Borders.Weight = -4138

# --- Move the active selection to I9 ---
$ws.Range("I9").Select()
